$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.058.22'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').Value = '2.010.34'
$ws.Range('E3').Value = '  -1.71%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '226.12'
$ws.Range('E5').Value = '  -1.27%  '
$ws.Range('D6').Value = '0.606'
$ws.Range('E6').Value = '  -1.12%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '55.07'
$ws.Range('E8').Value = '  -2.51%  '
$ws.Range('D9').Value = '0.374'
$ws.Range('E9').Value = '  -2.65%  '
$ws.Range('D10').Value = '0.0776'
$ws.Range('E10').Value = '  -4.10%  '
$ws.Range('E11').Value = '  -4.77%  '
$ws.Range('D12').Value = '2.308.00'
$ws.Range('E12').Value = '  -1.72%  '
$ws.Range('D13').Value = '14.01'
$ws.Range('E13').Value = '  -3.74%  '
$ws.Range('D14').Value = '19.74'
$ws.Range('E14').Value = '  -3.96%  '
$ws.Range('D15').Value = '5.19'
$ws.Range('E15').Value = '  -1.56%  '
$ws.Range('D16').Value = '0.735'
$ws.Range('E16').Value = '  -2.36%  '
$ws.Range('D17').Value = '2.008.13'
$ws.Range('E17').Value = '  -1.84%  '
$ws.Range('D18').Value = '37.020.52'
$ws.Range('E18').Value = '  -0.48%  '
$ws.Range('D19').Value = '6.20'
$ws.Range('E19').Value = '  +3.89%  '
$ws.Range('D20').Value = '68.30'
$ws.Range('E20').Value = '  -2.06%  '
$ws.Range('D21').Value = '0.0₃0810'
$ws.Range('E21').Value = '  -3.60%  '
$ws.Range('D22').Value = '223.67'
$ws.Range('E22').Value = '  -0.99%  '
$ws.Range('D24').Value = '2.43'
$ws.Range('E24').Value = '  +2.56%  '
$ws.Range('D25').Value = '2.16'
$ws.Range('E25').Value = '  -5.10%  '
$ws.Range('D26').Value = '164.30'
$ws.Range('E26').Value = '  -2.06%  '
$ws.Range('D27').Value = '8.92'
$ws.Range('E27').Value = '  -5.88%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '18.58'
$ws.Range('E28').Value = '  -1.57%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').Value = '0.124'
$ws.Range('E29').Value = '  -3.46%  '
$ws.Range('E30').Value = '  -7.37%  '
$ws.Range('E31').Value = '  -1.31%  '
$ws.Range('D32').Value = '4.44'
$ws.Range('E32').Value = '  -1.41%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '4.47'
$ws.Range('E33').Value = '  -1.98%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.0598'
$ws.Range('E34').Value = '  -2.08%  '
$ws.Range('D35').Value = '2.31'
$ws.Range('E35').Value = '  -3.29%  '
$ws.Range('E36').Value = '  +2.11%  '
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('D38').Value = '3.11'
$ws.Range('E38').Value = '  -2.30%  '
$ws.Range('D39').Value = '5.31'
$ws.Range('E39').Value = '  -1.55%  '
$ws.Range('D40').Value = '1.452.78'
$ws.Range('E40').Value = '  -1.62%  '
$ws.Range('D41').Value = '0.0211'
$ws.Range('E41').Value = '  -3.93%  '
$ws.Range('D42').Value = '94.49'
$ws.Range('E42').Value = '  -1.05%  '
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').Value = '0.0906'
$ws.Range('E43').Value = '  -3.53%  '
$ws.Range('B44').Value = 'HuobiToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D44').Value = '2.76'
$ws.Range('E44').Value = '  -4.85%  '
$ws.Range('D45').Value = '15.92'
$ws.Range('E45').Value = '  -6.02%  '
$ws.Range('E46').Value = '  -3.13%  '
$ws.Range('D47').Value = '7.11'
$ws.Range('E47').Value = '  +0.08%  '
$ws.Range('D48').Value = '0.994'
$ws.Range('E48').Value = '  -2.07%  '
$ws.Range('B49').Value = 'FTXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D49').Value = '3.85'
$ws.Range('E49').Value = '  +6.37%  '
$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D50').Value = '2.89'
$ws.Range('E50').Value = '  -0.55%  '
$ws.Range('D51').Value = '2.197.34'
$ws.Range('E51').Value = '  -1.69%  '
